$wb = $excel.ActiveWorkbook

# --- 1. Create the new "2022-Q4" sheet --------------------------------
# Duplicate the existing "2022-Q3" sheet (same column layout/styling as
# every other quarterly sheet) and drop the copy directly in front of it,
# i.e. right after "总计" — matching the new tab order:
#   总计, 2022-Q4, 2022-Q3, 2022-Q1, 2021-Q1, 2020-Q4
$templateSheet = $wb.Worksheets.Item("2022-Q3")
$templateSheet.Copy($templateSheet)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# Update the copied rows with the 2022-Q4 holdings figures (fund codes
# and names are unchanged quarter over quarter).
$newSheet.Range("D2").Value = "10.97"
$newSheet.Range("E2").Value = "94.25"
$newSheet.Range("F2").Value = "1.56"
$newSheet.Range("G2").Value = "0.1711"
$newSheet.Range("H2").Value = 6

$newSheet.Range("D3").Value = "3.29"
$newSheet.Range("E3").Value = "94.25"
$newSheet.Range("F3").Value = "1.56"
$newSheet.Range("G3").Value = "0.0513"
$newSheet.Range("H3").Value = 6

# --- 2. Insert a row for 2022-Q4 at the top of the "总计" summary sheet --
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.22

# Re-number the index column (A) for the rows that shifted down one slot.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
